$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right before the current row 364. This pushes the
# existing rows 364-409 down to 367-412 (carrying formatting, including the
# date-style on column D), and naturally yields the trailing three rows
# (410-412) that duplicate what used to be rows 407-409.
$ws.Rows("364:366").Insert()

# Row 364: new "Conconina(o)" observation dated 2021-11-05 (serial 44505)
$ws.Cells.Item(364, 1).Value = 11
$ws.Cells.Item(364, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(364, 3).Value = "Bíobío"
$ws.Cells.Item(364, 4).Value = 44505
$ws.Cells.Item(364, 5).Value = 8
$ws.Cells.Item(364, 6).Value = 100112033
$ws.Cells.Item(364, 7).Value = "Lechuga"
$ws.Cells.Item(364, 8).Value = "Conconina(o)"
$ws.Cells.Item(364, 9).Value = "Primera"
$ws.Cells.Item(364, 10).Value = 250
$ws.Cells.Item(364, 11).Value = 4000
$ws.Cells.Item(364, 12).Value = 4500
$ws.Cells.Item(364, 13).Value = 4300
$ws.Cells.Item(364, 14).Value = "$/caja 10 unidades"
$ws.Cells.Item(364, 15).Value = "Región Metropolitana"
$ws.Cells.Item(364, 16).Value = 430
$ws.Cells.Item(364, 17).Value = 10
$ws.Cells.Item(364, 18).Value = "Hortaliza"

# Row 365: new "Escarola" observation dated 2021-11-05 (serial 44505)
$ws.Cells.Item(365, 1).Value = 11
$ws.Cells.Item(365, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(365, 3).Value = "Bíobío"
$ws.Cells.Item(365, 4).Value = 44505
$ws.Cells.Item(365, 5).Value = 8
$ws.Cells.Item(365, 6).Value = 100112033
$ws.Cells.Item(365, 7).Value = "Lechuga"
$ws.Cells.Item(365, 8).Value = "Escarola"
$ws.Cells.Item(365, 9).Value = "Primera"
$ws.Cells.Item(365, 10).Value = 220
$ws.Cells.Item(365, 11).Value = 5500
$ws.Cells.Item(365, 12).Value = 6000
$ws.Cells.Item(365, 13).Value = 5727
$ws.Cells.Item(365, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(365, 15).Value = "Región del Maule"
$ws.Cells.Item(365, 16).Value = 382
$ws.Cells.Item(365, 17).Value = 15
$ws.Cells.Item(365, 18).Value = "Hortaliza"

# Row 366: new "Marina" observation dated 2021-11-05 (serial 44505)
$ws.Cells.Item(366, 1).Value = 11
$ws.Cells.Item(366, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(366, 3).Value = "Bíobío"
$ws.Cells.Item(366, 4).Value = 44505
$ws.Cells.Item(366, 5).Value = 8
$ws.Cells.Item(366, 6).Value = 100112033
$ws.Cells.Item(366, 7).Value = "Lechuga"
$ws.Cells.Item(366, 8).Value = "Marina"
$ws.Cells.Item(366, 9).Value = "Primera"
$ws.Cells.Item(366, 10).Value = 170
$ws.Cells.Item(366, 11).Value = 4000
$ws.Cells.Item(366, 12).Value = 4500
$ws.Cells.Item(366, 13).Value = 4265
$ws.Cells.Item(366, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(366, 15).Value = "Región del Maule"
$ws.Cells.Item(366, 16).Value = 237
$ws.Cells.Item(366, 17).Value = 18
$ws.Cells.Item(366, 18).Value = "Hortaliza"
